# ccb_cause_to_100_year_cause.xlsx
# Reclassify D02 ("Endocrine, blood, immune disorders") and D03 ("Sickle
# cell disorders and trait") from "Cardiovascular" / "Other Cardiovascular"
# to "Other Chronic" across the OurGroup0 / OurGroup1 / OurGroup3 columns
# (rows 37 and 38 of the ccb_cause_to_100_year_causes sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C37:E37").Value = "Other Chronic"
$ws.Range("C38:E38").Value = "Other Chronic"
